# activitiesClassification.xlsx - "divisao de tarefas + definicao de recordes"
#
# Adds a task-assignment column (K) naming who is responsible for each
# activity row (cycling through Oliveira / Santos / Camposinhos), marks the
# previously-blank checkbox in J61 with an "x" like its sibling cells, and
# moves the view/selection to where the work was happening (around I55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K21:K61 - round-robin task owner for each activity row ---------------
$names = @("Oliveira", "Santos", "Camposinhos")
for ($r = 21; $r -le 61; $r++) {
    $owner = $names[($r - 21) % 3]
    $ws.Cells.Item($r, 11).Value2 = $owner
}

# --- J61 - complete the missing "x" mark (and drop its stray underline) ---
$j61 = $ws.Range("J61")
$j61.Font.Underline = -4142   # xlUnderlineStyleNone
$j61.Value2 = "x"

# --- move the view to where the edits were made ---------------------------
$ws.Activate() | Out-Null
$ws.Range("I55").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 2
